# Scheduled refresh of currentAveragePrice / LevePrice / LeveProfit columns
# (H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
#  K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28: The Writing Is Not on the Wall / Enchanted Silver Ink
$ws.Range("H28").Value = 1856.1428
$ws.Range("I28").Value = 1798.8
$ws.Range("J28").Value = 1999.5
$ws.Range("K28").Value = 1798.8
$ws.Range("L28").Value = 1999.5
$ws.Range("M28").Value = -1313.8
$ws.Range("N28").Value = -2969.5
# Row 92: Whinier than the Sword / Enchanted Koppranickel Ink
$ws.Range("H92").Value = 400
$ws.Range("I92").Value = 400
$ws.Range("K92").Value = 400
$ws.Range("M92").Value = 848
# Row 103: Let Loose the Juice / Persimmon Tannin
$ws.Range("H103").Value = 3992.7144
$ws.Range("J103").Value = 3990.8333
$ws.Range("L103").Value = 11972.4999
$ws.Range("N103").Value = -13144.4999

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 3599.8
$ws.Range("I86").Value = 2666.3333
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 2666.3333
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -1543.3333
$ws.Range("N86").Value = -7246
# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 3599.8
$ws.Range("I89").Value = 2666.3333
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 13331.6665
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -7715.666499999999
$ws.Range("N89").Value = -36232
# Row 94: High Steal / High Steel Nugget
$ws.Range("H94").Value = 914.5
$ws.Range("I94").Value = 914.5
$ws.Range("K94").Value = 914.5
$ws.Range("M94").Value = -463.5
# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 1832.6666
$ws.Range("J99").Value = 1000
$ws.Range("L99").Value = 1000
$ws.Range("N99").Value = -3996
# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 2124.75
$ws.Range("I105").Value = 1166.3334
$ws.Range("K105").Value = 1166.3334
$ws.Range("M105").Value = 580.6666

$ws = $wb.Worksheets.Item("CRP")
# Row 43: The Long Lance of the Law / Steel Halberd
$ws.Range("H43").Value = 11839.4
$ws.Range("J43").Value = 11839.4
$ws.Range("L43").Value = 11839.4
$ws.Range("N43").Value = -12207.4
# Row 101: Everybody's Heard about the 'Berd / Doman Steel Halberd
$ws.Range("H101").Value = 11839.4
$ws.Range("J101").Value = 11839.4
$ws.Range("L101").Value = 11839.4
$ws.Range("N101").Value = -18329.4
# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 2149.7
$ws.Range("I132").Value = 1833
$ws.Range("K132").Value = 5499
$ws.Range("M132").Value = -2969
# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 3813.6843
$ws.Range("I134").Value = 2439.75
$ws.Range("K134").Value = 7319.25
$ws.Range("M134").Value = -4784.25

$ws = $wb.Worksheets.Item("CUL")
# Row 76: Old Victories, New Tastes / Dhalmel Fricassee
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").Value = $null
# Row 79: The Eats of Authenticity (L) / Dhalmel Fricassee
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").Value = $null
# Row 113: Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 313.83334
$ws.Range("J113").Value = 277.5
$ws.Range("L113").Value = 832.5
$ws.Range("N113").Value = -5172.5
# Row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 898.75
$ws.Range("I122").Value = 899
$ws.Range("K122").Value = 8091
$ws.Range("M122").Value = -5641

$ws = $wb.Worksheets.Item("GSM")
# Row 29: Music to Their Ears / Brass Ear Cuffs
$ws.Range("H29").Value = 1000
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = $null
# Row 42: It's Only Love / Silver Choker
$ws.Range("H42").Value = 14150
$ws.Range("I42").Value = 10000
$ws.Range("J42").Value = 18300
$ws.Range("K42").Value = 10000
$ws.Range("L42").Value = 18300
$ws.Range("M42").Value = -9515
$ws.Range("N42").Value = -19270
# Row 64: Halonic Hermeneutics / Yeti Staff
$ws.Range("H64").Value = 35000
$ws.Range("J64").Value = 35000
$ws.Range("L64").Value = 35000
$ws.Range("N64").Value = -35496
# Row 67: Transposing Theology (L) / Yeti Staff
$ws.Range("H67").Value = 35000
$ws.Range("J67").Value = 35000
$ws.Range("L67").Value = 35000
$ws.Range("N67").Value = -36716
# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 2495
$ws.Range("I80").Value = 2495
$ws.Range("J80").Value = 2495
$ws.Range("K80").Value = 2495
$ws.Range("L80").Value = 2495
$ws.Range("M80").Value = -1497
$ws.Range("N80").Value = -4491
# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 2495
$ws.Range("I83").Value = 2495
$ws.Range("J83").Value = 2495
$ws.Range("K83").Value = 12475
$ws.Range("L83").Value = 12475
$ws.Range("M83").Value = -7483
$ws.Range("N83").Value = -22459
# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Range("H113").Value = 1149.6666
$ws.Range("I113").Value = 1149.6666
$ws.Range("K113").Value = 1149.6666
$ws.Range("M113").Value = 1020.3334
# Row 115: Unsung Generosity / Manasilver Choker
$ws.Range("H115").Value = 14150
$ws.Range("I115").Value = 10000
$ws.Range("J115").Value = 18300
$ws.Range("K115").Value = 10000
$ws.Range("L115").Value = 18300
$ws.Range("M115").Value = -8825
$ws.Range("N115").Value = -20650
# Row 118: A Magnanimous Refrain / Triplite Earrings of Casting
$ws.Range("H118").Value = 9484.5
$ws.Range("J118").Value = 9484.5
$ws.Range("L118").Value = 9484.5
$ws.Range("N118").Value = -12798.5

$ws = $wb.Worksheets.Item("LTW")
# Row 12: A Place to Call Helm / Hard Leather Pot Helm
$ws.Range("H12").Value = 3459.25
$ws.Range("J12").Value = 1491.3334
$ws.Range("L12").Value = 1491.3334
$ws.Range("N12").Value = -1831.3334
# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 1250
$ws.Range("I61").Value = 500
$ws.Range("K61").Value = 500
$ws.Range("M61").Value = -298
# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 3250
$ws.Range("I82").Value = 3000
$ws.Range("K82").Value = 3000
$ws.Range("M82").Value = -2639
# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 3250
$ws.Range("I85").Value = 3000
$ws.Range("K85").Value = 3000
$ws.Range("M85").Value = -1752
# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 1250
$ws.Range("I113").Value = 500
$ws.Range("K113").Value = 500
$ws.Range("M113").Value = 1670
# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 6900
$ws.Range("I132").Value = 7500
$ws.Range("K132").Value = 22500
$ws.Range("M132").Value = -19970

$ws = $wb.Worksheets.Item("WVR")
# Row 14: Hat in Hand / Straw Hat
$ws.Range("H14").Value = 2500
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 2500
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = $null
$ws.Range("N14").Value = -2836
# Row 107: Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value = 802
$ws.Range("I107").Value = 802
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2406
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -486
$ws.Range("N107").Value = $null
# Row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 1896.75
$ws.Range("I113").Value = 1896.75
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 5690.25
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -3520.25
$ws.Range("N113").Value = $null
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 4477.8887
$ws.Range("I136").Value = 4477.8887
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 13433.6661
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -10883.6661
$ws.Range("N136").Value = $null
